$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = "648 74 65 62"
$ws.Range("N2").Value = "2025-09-17 13:14:08"
$ws.Range("O2").Value = "2025-09-17 13:14:08"

# Row 3
$ws.Range("J3").Value = "635 57 00 94"
$ws.Range("N3").Value = "2025-09-17 13:25:58"
$ws.Range("O3").Value = "2025-09-17 13:25:58"

# Row 4
$ws.Range("J4").Value = "630 23 45 85"
$ws.Range("N4").Value = "2025-09-17 13:25:58"
$ws.Range("O4").Value = "2025-09-17 13:25:58"

# Row 5
$ws.Range("J5").Value = "644 49 22 12"
$ws.Range("N5").Value = "2025-09-17 13:25:58"
$ws.Range("O5").Value = "2025-09-17 13:25:58"

# Row 6
$ws.Range("J6").Value = "644 49 22 12"
$ws.Range("N6").Value = "2025-09-17 13:25:58"
$ws.Range("O6").Value = "2025-09-17 13:25:58"
